$d = $word.ActiveDocument

# Locate the paragraph about the effect of nebulosity using a unique
# anchor phrase that spans the text to be trimmed.
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute(
    "pela superfície do planeta gera",  # FindText
    $false,                              # MatchCase
    $false,                              # MatchWholeWord
    $false,                              # MatchWildcards
    $false,                              # MatchSoundsLike
    $false,                              # MatchAllWordForms
    $true,                               # Forward
    1,                                   # Wrap (wdFindContinue)
    $false,                              # Format
    $null,                               # ReplaceWith
    0                                    # Replace (wdReplaceNone)
)

if (-not $found) {
    throw "Could not find the target sentence about nebulosity."
}

$matchRange = $find.Parent
$matchStart = $matchRange.Start
$matchText = $matchRange.Text

# Remove "do planeta " so "... emitida pela superfície " is immediately
# followed by "gera um aquecimento. ..."
$needle = "do planeta "
$relIdx = $matchText.IndexOf($needle)
$cutStart = $matchStart + $relIdx
$cutEnd = $cutStart + $needle.Length

$cutRange = $d.Range($cutStart, $cutEnd)
$cutRange.Text = ""

# Re-seat the _GoBack bookmark (last-edit marker) exactly at the cut
# point; Word splits the underlying run into two at that text position
# when the bookmark is (re)inserted there, matching how the paragraph
# was actually edited.
$splitPoint = $d.Range($cutStart, $cutStart)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null
